# Update Louisiana overview factsheet per COMM text edits.
#
# The underlying data values are unchanged; the edit re-types a batch of
# numeric "count" cells as literal text (so they round-trip/display exactly
# as authored, e.g. "1,214" with a thousands separator) and refreshes the
# "St. Helena Parish" placeholder row (was all zeros) plus adds a new
# "Total" row to the County sheet.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$addr, [string]$text)
    $rng = $ws.Range($addr)
    # Force text storage so a numeric-looking string (e.g. "7", "1,214")
    # is kept as text rather than being re-parsed back into a number.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    # Drop the now-unnecessary "Text" number-format override so the cell
    # doesn't carry a style it didn't have before.
    $rng.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet "Overall": A2 1214 -> "1,214"
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall "A2" "1,214"

# ---------------------------------------------------------------------
# Sheet "County": B2:B64 numbers -> text (same digits), row 65
# (St. Helena Parish) refreshed, and new "Total" row 66 appended.
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

$countyCounts = @{
    2 = "7";   3 = "4";   4 = "6";   5 = "4";   6 = "7";   7 = "6";
    8 = "4";   9 = "13";  10 = "81"; 11 = "35"; 12 = "1";  13 = "1";
    14 = "1";  15 = "2";  16 = "7";  17 = "4";  18 = "177";19 = "2";
    20 = "1";  21 = "7";  22 = "4";  23 = "1";  24 = "8";  25 = "6";
    26 = "3";  27 = "6";  28 = "80"; 29 = "2";  30 = "97"; 31 = "20";
    32 = "17"; 33 = "7";  34 = "4";  35 = "8";  36 = "13"; 37 = "305";
    38 = "52"; 39 = "2";  40 = "3";  41 = "36"; 42 = "2";  43 = "2";
    44 = "2";  45 = "4";  46 = "6";  47 = "2";  48 = "3";  49 = "15";
    50 = "4";  51 = "10"; 52 = "38"; 53 = "26"; 54 = "1";  55 = "22";
    56 = "6";  57 = "2";  58 = "2";  59 = "10"; 60 = "5";  61 = "2";
    62 = "2";  63 = "2";  64 = "2"
}

foreach ($r in $countyCounts.Keys) {
    $addr = "B" + $r
    $val = $countyCounts[$r]
    Set-TextValue $wsCounty $addr $val
}

# Row 65, "St. Helena Parish" — was a placeholder row of zeros.
Set-TextValue $wsCounty "B65" "0.00%"
Set-TextValue $wsCounty "C65" "`$0"
Set-TextValue $wsCounty "D65" "0.00%"
Set-TextValue $wsCounty "E65" "0.00%"
Set-TextValue $wsCounty "F65" "0.00%"

# New row 66 — statewide "Total" (mirrors the Total row already present
# on the other sheets).
Set-TextValue $wsCounty "A66" "Total"
Set-TextValue $wsCounty "B66" "1,214"
Set-TextValue $wsCounty "C66" "`$2,869,770,671"
Set-TextValue $wsCounty "D66" "7.97%"
Set-TextValue $wsCounty "E66" "-27.48%"
Set-TextValue $wsCounty "F66" "72.08%"

# ---------------------------------------------------------------------
# Sheet "Congressional District": B2:B7 numbers -> text, B8 Total -> "1,214"
# ---------------------------------------------------------------------
$wsCD = $wb.Worksheets.Item("Congressional District")

$cdCounts = @{ 2 = "149"; 3 = "362"; 4 = "168"; 5 = "167"; 6 = "188"; 7 = "180" }
foreach ($r in $cdCounts.Keys) {
    $addr = "B" + $r
    $val = $cdCounts[$r]
    Set-TextValue $wsCD $addr $val
}
Set-TextValue $wsCD "B8" "1,214"

# ---------------------------------------------------------------------
# Sheet "Size": B2:B7 numbers -> text, B8 Total -> "1,214"
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")

$sizeCounts = @{ 2 = "372"; 3 = "330"; 4 = "195"; 5 = "99"; 6 = "165"; 7 = "53" }
foreach ($r in $sizeCounts.Keys) {
    $addr = "B" + $r
    $val = $sizeCounts[$r]
    Set-TextValue $wsSize $addr $val
}
Set-TextValue $wsSize "B8" "1,214"

# ---------------------------------------------------------------------
# Sheet "Subsector": B2:B12 numbers -> text, B13 Total -> "1,214"
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")

$subCounts = @{
    2 = "106"; 3 = "189"; 4 = "32"; 5 = "106"; 6 = "17"; 7 = "397";
    8 = "3";   9 = "116"; 10 = "19"; 11 = "221"; 12 = "8"
}
foreach ($r in $subCounts.Keys) {
    $addr = "B" + $r
    $val = $subCounts[$r]
    Set-TextValue $wsSub $addr $val
}
Set-TextValue $wsSub "B13" "1,214"
